$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that no longer exist after recomputation (old rows 10 and 11)
$ws.Range("A10:A11").EntireRow.Delete()

# Row 2
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Ebi3"
$ws.Range("C2").Value = "Il27ra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.03073166666667
$ws.Range("H2").Value = 30.092195
$ws.Range("I2").Value = 0.5703560915538793
$ws.Range("J2").Value = 0.5703560915538795
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.171010666666667
$ws.Range("N2").Value = 6.513032
$ws.Range("O2").Value = 0.5189085589295733
$ws.Range("P2").Value = 0.5189085589295733
$ws.Range("Q2").Value = 21.77682544280444
$ws.Range("R2").Value = 195.99142898524
$ws.Range("S2").Value = 0.2959626575449273
$ws.Range("T2").Value = 0.2959626575449274

# Row 3
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Ebi3"
$ws.Range("C3").Value = "Il27ra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.03073166666667
$ws.Range("H3").Value = 30.092195
$ws.Range("I3").Value = 0.5703560915538793
$ws.Range("J3").Value = 0.5703560915538795
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.561521333333333
$ws.Range("N3").Value = 4.684564
$ws.Range("O3").Value = 0.3732302181922886
$ws.Range("P3").Value = 0.3732302181922886
$ws.Range("Q3").Value = 15.66320148644222
$ws.Range("R3").Value = 140.96881337798
$ws.Range("S3").Value = 0.2128741284979553
$ws.Range("T3").Value = 0.2128741284979553

# Row 4
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Ebi3"
$ws.Range("C4").Value = "Il27ra"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.03073166666667
$ws.Range("H4").Value = 30.092195
$ws.Range("I4").Value = 0.5703560915538793
$ws.Range("J4").Value = 0.5703560915538795
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4199566666666667
$ws.Range("N4").Value = 1.25987
$ws.Range("O4").Value = 0.1003768024076347
$ws.Range("P4").Value = 0.1003768024076346
$ws.Range("Q4").Value = 4.212472634961111
$ws.Range("R4").Value = 37.91225371465
$ws.Range("S4").Value = 0.05725052070389453
$ws.Range("T4").Value = 0.05725052070389453

# Row 5
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Ebi3"
$ws.Range("C5").Value = "Il27ra"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.03073166666667
$ws.Range("H5").Value = 30.092195
$ws.Range("I5").Value = 0.5703560915538793
$ws.Range("J5").Value = 0.5703560915538795
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03131333333333333
$ws.Range("N5").Value = 0.09394
$ws.Range("O5").Value = 0.007484420470503464
$ws.Range("P5").Value = 0.007484420470503464
$ws.Range("Q5").Value = 0.3140956442555555
$ws.Range("R5").Value = 2.8268607983
$ws.Range("S5").Value = 0.004268784807102202
$ws.Range("T5").Value = 0.004268784807102203

# Row 6
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Ebi3"
$ws.Range("C6").Value = "Il27ra"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 7.556056333333333
$ws.Range("H6").Value = 22.668169
$ws.Range("I6").Value = 0.4296439084461207
$ws.Range("J6").Value = 0.4296439084461207
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.171010666666667
$ws.Range("N6").Value = 6.513032
$ws.Range("O6").Value = 0.5189085589295733
$ws.Range("P6").Value = 0.5189085589295733
$ws.Range("Q6").Value = 16.40427889760089
$ws.Range("R6").Value = 147.638510078408
$ws.Range("S6").Value = 0.222945901384646
$ws.Range("T6").Value = 0.222945901384646

# Row 7
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Ebi3"
$ws.Range("C7").Value = "Il27ra"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 7.556056333333333
$ws.Range("H7").Value = 22.668169
$ws.Range("I7").Value = 0.4296439084461207
$ws.Range("J7").Value = 0.4296439084461207
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.561521333333333
$ws.Range("N7").Value = 4.684564
$ws.Range("O7").Value = 0.3732302181922886
$ws.Range("P7").Value = 0.3732302181922886
$ws.Range("Q7").Value = 11.79894316036844
$ws.Range("R7").Value = 106.190488443316
$ws.Range("S7").Value = 0.1603560896943333
$ws.Range("T7").Value = 0.1603560896943333

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Ebi3"
$ws.Range("C8").Value = "Il27ra"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.556056333333333
$ws.Range("H8").Value = 22.668169
$ws.Range("I8").Value = 0.4296439084461207
$ws.Range("J8").Value = 0.4296439084461207
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4199566666666667
$ws.Range("N8").Value = 1.25987
$ws.Range("O8").Value = 0.1003768024076347
$ws.Range("P8").Value = 0.1003768024076346
$ws.Range("Q8").Value = 3.173216230892222
$ws.Range("R8").Value = 28.55894607803
$ws.Range("S8").Value = 0.04312628170374013
$ws.Range("T8").Value = 0.04312628170374012

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Ebi3"
$ws.Range("C9").Value = "Il27ra"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.556056333333333
$ws.Range("H9").Value = 22.668169
$ws.Range("I9").Value = 0.4296439084461207
$ws.Range("J9").Value = 0.4296439084461207
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.03131333333333333
$ws.Range("N9").Value = 0.09394
$ws.Range("O9").Value = 0.007484420470503464
$ws.Range("P9").Value = 0.007484420470503464
$ws.Range("Q9").Value = 0.2366053106511111
$ws.Range("R9").Value = 2.12944779586
$ws.Range("S9").Value = 0.003215635663401262
$ws.Range("T9").Value = 0.003215635663401262

